$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (columns A..AG), per updated spot price data
$ws.Range("A2").Value = 46058
$ws.Range("B2").Value = 2.16
$ws.Range("C2").Value = 0.03
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = -0.08
$ws.Range("F2").Value = -0.1
$ws.Range("G2").Value = -0.1
$ws.Range("H2").Value = -0.1
$ws.Range("I2").Value = -0.01
$ws.Range("J2").Value = 0.59
$ws.Range("K2").Value = 2.44
$ws.Range("L2").Value = 2.21
$ws.Range("M2").Value = 0.24
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = -0.04
$ws.Range("Q2").Value = -0.1
$ws.Range("R2").Value = -0.02
$ws.Range("S2").Value = 0.34
$ws.Range("T2").Value = 1.41
$ws.Range("U2").Value = 9.07
$ws.Range("V2").Value = 11.79
$ws.Range("W2").Value = 11.7
$ws.Range("X2").Value = 6.16
$ws.Range("Y2").Value = 0.75
$ws.Range("Z2").Value = 2.01
$ws.Range("AB2").Value = 7.6
$ws.Range("AD2").Value = 11.74
$ws.Range("AF2").Value = 5.24
$ws.Range("AG2").Value = "1h-23h"
